$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "updated activity till excel form" - refresh the batting-innings rows for
# Kedar Jadhav (Chennai Super Kings): the existing 4 innings rows get new
# figures, and the innings that used to sit in row 2 is pushed down to a
# brand-new row 6 (sheet grows from A1:F5 to A1:F6).

$nbsp = [char]0x00A0
$playerName = "Kedar Jadhav$nbsp"
$teamName = "Chennai Super Kings"

# Keep the numeric-looking stats stored as TEXT, matching how this sheet
# already stores every value (runs/balls/fours/sixes as text strings).
$ws.Range("C2:F6").NumberFormat = "@"

# Row 2
$ws.Range("C2").Value = "7"
$ws.Range("D2").Value = "12"
$ws.Range("E2").Value = "1"
$ws.Range("F2").Value = "0"

# Row 3
$ws.Range("C3").Value = "22"
$ws.Range("D3").Value = "16"
$ws.Range("E3").Value = "3"
$ws.Range("F3").Value = "0"

# Row 4
$ws.Range("C4").Value = "4"
$ws.Range("D4").Value = "7"
$ws.Range("E4").Value = "0"
$ws.Range("F4").Value = "0"

# Row 5
$ws.Range("C5").Value = "3"
$ws.Range("D5").Value = "10"
$ws.Range("E5").Value = "0"
$ws.Range("F5").Value = "0"

# New row 6 - carries forward what used to be row 2's figures
$ws.Range("A6").Value = $playerName
$ws.Range("B6").Value = $teamName
$ws.Range("C6").Value = "26"
$ws.Range("D6").Value = "21"
$ws.Range("E6").Value = "3"
$ws.Range("F6").Value = "0"
